$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $text)
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

$ws.Range("D2").Value = "69.128.99"
$ws.Range("E2").Value = "  +3.51%  "
$ws.Range("D3").Value = "3.746.47"
$ws.Range("E3").Value = "  +1.44%  "
$ws.Range("E4").Value = "  -0.57%  "
Set-TextValue $ws.Range("D5") "601.72"
$ws.Range("E5").Value = "  +3.18%  "
Set-TextValue $ws.Range("D6") "168.34"
$ws.Range("E6").Value = "  +0.33%  "
$ws.Range("D7").Value = "3.743.35"
$ws.Range("E7").Value = "  +1.68%  "
$ws.Range("E8").Value = "  -0.27%  "
$ws.Range("E9").Value = "  +4.58%  "
Set-TextValue $ws.Range("D10") "0.166"
$ws.Range("E10").Value = "  +8.20%  "
Set-TextValue $ws.Range("D11") "6.32"
$ws.Range("E11").Value = "  +4.70%  "
Set-TextValue $ws.Range("D12") "0.461"
$ws.Range("E12").Value = "  +2.82%  "
Set-TextValue $ws.Range("D13") "38.43"
$ws.Range("E13").Value = "  +4.52%  "
Set-TextValue $ws.Range("D14") "0.0000245"
$ws.Range("E14").Value = "  +4.00%  "
$ws.Range("D15").Value = "4.370.37"
$ws.Range("E15").Value = "  -0.90%  "
$ws.Range("D16").Value = "3.742.16"
$ws.Range("E16").Value = "  -0.65%  "
$ws.Range("D17").Value = "69.071.86"
$ws.Range("E17").Value = "  +3.21%  "
Set-TextValue $ws.Range("D18") "7.28"
$ws.Range("E18").Value = "  +4.65%  "
$ws.Range("E19").Value = "  +1.17%  "
Set-TextValue $ws.Range("D20") "17.32"
$ws.Range("E20").Value = "  +11.09%  "
Set-TextValue $ws.Range("D21") "499.64"
$ws.Range("E21").Value = "  +5.20%  "
Set-TextValue $ws.Range("D22") "10.30"
$ws.Range("E22").Value = "  +16.96%  "
Set-TextValue $ws.Range("D23") "0.726"
$ws.Range("E23").Value = "  +3.47%  "
Set-TextValue $ws.Range("D24") "85.28"
$ws.Range("E24").Value = "  +3.66%  "
$ws.Range("E25").Value = "  +6.80%  "
$ws.Range("E26").Value = "  +1.65%  "
Set-TextValue $ws.Range("D27") "12.31"
$ws.Range("E27").Value = "  +3.56%  "
Set-TextValue $ws.Range("D28") "10.14"
$ws.Range("E28").Value = "  +2.46%  "
$ws.Range("E29").Value = "  +0.58%  "
$ws.Range("E30").Value = "  +3.71%  "
$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D31") "7.97"
$ws.Range("E31").Value = "  +6.51%  "
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D32") "2.42"
$ws.Range("E32").Value = "  +3.93%  "
Set-TextValue $ws.Range("D33") "31.83"
$ws.Range("E33").Value = "  +1.18%  "
$ws.Range("D34").Value = "3.884.79"
$ws.Range("E34").Value = "  +0.35%  "
$ws.Range("E35").Value = "  +3.60%  "
$ws.Range("D36").Value = "3.675.33"
$ws.Range("E36").Value = "  +0.30%  "
$ws.Range("E37").Value = "  -0.79%  "
$ws.Range("E38").Value = "  +3.59%  "
Set-TextValue $ws.Range("D39") "5.81"
$ws.Range("E39").Value = "  +4.79%  "
$ws.Range("E40").Value = "  +2.01%  "
Set-TextValue $ws.Range("D41") "0.325"
$ws.Range("E41").Value = "  +3.58%  "
Set-TextValue $ws.Range("D42") "441.45"
$ws.Range("E42").Value = "  +0.46%  "
Set-TextValue $ws.Range("D43") "49.03"
$ws.Range("E43").Value = "  +1.32%  "
Set-TextValue $ws.Range("D44") "2.00"
$ws.Range("E44").Value = "  +2.94%  "
$ws.Range("E45").Value = "  +5.03%  "
Set-TextValue $ws.Range("D46") "8.42"
$ws.Range("E46").Value = "  +4.47%  "
$ws.Range("E47").Value = "  +0.01%  "
Set-TextValue $ws.Range("D48") "40.58"
$ws.Range("E48").Value = "  +0.74%  "
Set-TextValue $ws.Range("D49") "143.07"
$ws.Range("E49").Value = "  +2.22%  "
$ws.Range("E50").Value = "  +3.95%  "
$ws.Range("D51").Value = "2.750.15"
$ws.Range("E51").Value = "  -0.55%  "
